$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.592.33'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '1.966.80'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('D4').Value = "'1.013"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').Value = "'323.40"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = "'0.4803"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -3.61%  '
$ws.Range('D8').Value = "'0.4069"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('D9').Value = "'54.15"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = "'0.08522"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.45%  '
$ws.Range('D11').Value = "'1.064"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('D12').Value = "'22.47"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('D13').Value = '1.998.20'
$ws.Range('E13').Value = '  +4.36%  '
$ws.Range('D14').Value = "'7.639"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').Value = "'6.203"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.67%  '
$ws.Range('D16').Value = "'1.014"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = "'91.29"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').Value = "'0.06646"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = "'18.64"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.70%  '
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').Value = "'5.875"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').Value = '28.621.04'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').Value = "'11.53"
$ws.Range('D24').ClearFormats()
$ws.Range('D25').Value = "'2.303"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').Value = '2.230.46'
$ws.Range('E26').Value = '  +3.02%  '
$ws.Range('D27').Value = "'156.56"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('D28').Value = "'20.40"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').Value = "'5.911"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.03%  '
$ws.Range('D30').Value = "'2.185"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('D32').Value = "'0.9926"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.55%  '
$ws.Range('D33').Value = "'0.09683"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('E34').Value = '  -4.38%  '
$ws.Range('D35').Value = "'3.712"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = "'5.650"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').Value = "'9.162"
$ws.Range('D37').ClearFormats()
$ws.Range('D38').Value = "'0.02341"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.04%  '
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('D41').Value = "'0.6247"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.86%  '
$ws.Range('D42').Value = "'11.24"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('D44').Value = "'0.1926"
$ws.Range('D44').ClearFormats()
$ws.Range('D45').Value = "'1.364"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.41%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.5978"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.31%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'13.07"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').Value = "'2.075"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.20%  '
$ws.Range('D49').Value = "'3.418"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('D50').Value = "'0.06838"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').Value = "'111.66"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.80%  '
